# Weekly fruit/vegetable price update: a new daily price record is inserted
# as row 66 on the active sheet, pushing the existing rows 66-82 down to
# rows 67-83 (dimension grows from A1:R82 to A1:R83).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 66 (shifts rows 66..82 -> 67..83).
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new record.
$ws.Range("A66").Value = 3
$ws.Range("B66").Value = "Femacal de La Calera"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 45218
$ws.Range("E66").Value = 5
$ws.Range("F66").Value = 300000000
$ws.Range("G66").Value = "Espárragos"
$ws.Range("H66").Value = "Verde"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 2080
$ws.Range("K66").Value = 1700
$ws.Range("L66").Value = 1800
$ws.Range("M66").Value = 1747
$ws.Range("N66").Value = "$/kilo"
$ws.Range("O66").Value = "Provincia de Linares"
$ws.Range("P66").Value = 1747
$ws.Range("Q66").Value = 1
$ws.Range("R66").Value = "Hortaliza"
